$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Node" column (header "Node" / value placeholder "[[Data.Node.Name]]")
# is being removed from the table in row 6/7. DeliveryRoutes and
# DeliveryTrips shift one column to the left (E, F) and the former
# column G is cleared out entirely.

$ws.Range("E6").Value2 = $ws.Range("F6").Value2
$ws.Range("F6").Value2 = $ws.Range("G6").Value2
$ws.Range("G6").Clear()

$ws.Range("E7").Value2 = $ws.Range("F7").Value2
$ws.Range("F7").Value2 = $ws.Range("G7").Value2
$ws.Range("G7").Clear()
